$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns P and Q need the same formatting as column O (bold/bordered header style)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Header row additions
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update existing columns I, K, M, O (swap 1s and 2s) and add P, Q columns for data rows 2-25
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I
    $ws.Cells.Item($r, 11).Value = 1  # K
    $ws.Cells.Item($r, 13).Value = 2  # M
    $ws.Cells.Item($r, 15).Value = 1  # O
    $ws.Cells.Item($r, 16).Value = 2  # P
    $ws.Cells.Item($r, 17).Value = 2  # Q
}
